# Update the crypto "Price" (D) and "Volume(1h)" (E) columns with the
# latest scraped symbol-list snapshot. All of these cells hold text
# values (e.g. "306.68", "-2.62%"), not numbers, so every assignment is
# prefixed with a literal apostrophe to force Excel to store it as text
# instead of re-interpreting it as a number/percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.68"

$ws.Range("D3").Value = "'40.56"
$ws.Range("E3").Value = "'-2.62%"

$ws.Range("D4").Value = "'5.048"
$ws.Range("E4").Value = "'-2.83%"

$ws.Range("D5").Value = "'0.07605"
$ws.Range("E5").Value = "'-6.04%"

$ws.Range("D6").Value = "'4.242"
$ws.Range("E6").Value = "'-2.74%"

$ws.Range("E7").Value = "'-9.06%"

$ws.Range("D8").Value = "'0.9045"
$ws.Range("E8").Value = "'-2.69%"

$ws.Range("E9").Value = "'-12.16%"

$ws.Range("D10").Value = "'0.1757"
$ws.Range("E10").Value = "'-5.20%"

$ws.Range("D11").Value = "'0.09184"
$ws.Range("E11").Value = "'-0.91%"

$ws.Range("D12").Value = "'0.04347"
$ws.Range("E12").Value = "'-5.17%"

$ws.Range("E13").Value = "'-0.21%"

$ws.Range("D14").Value = "'0.001258"
$ws.Range("E14").Value = "'-2.78%"

$ws.Range("D15").Value = "'0.005817"
$ws.Range("E15").Value = "'-0.45%"

$ws.Range("E16").Value = "'0.72%"

$ws.Range("E17").Value = "'-6.36%"

$ws.Range("E18").Value = "'-3.10%"

$ws.Range("D19").Value = "'6.824"
$ws.Range("E19").Value = "'-7.86%"

$ws.Range("D20").Value = "'0.1350"
$ws.Range("E20").Value = "'-2.38%"

$ws.Range("D21").Value = "'0.2847"
$ws.Range("E21").Value = "'11.55%"

$ws.Range("D22").Value = "'0.04163"
$ws.Range("E22").Value = "'-0.58%"

$ws.Range("D23").Value = "'0.001216"
$ws.Range("E23").Value = "'-2.28%"

$ws.Range("D24").Value = "'0.004064"
$ws.Range("E24").Value = "'-4.26%"

$ws.Range("D25").Value = "'0.0001303"
$ws.Range("E25").Value = "'6.37%"

$ws.Range("D38").Value = "'0.02410"
$ws.Range("E38").Value = "'-6.71%"

$ws.Range("D39").Value = "'0.05129"
$ws.Range("E39").Value = "'-6.32%"

$ws.Range("D40").Value = "'0.007854"
$ws.Range("E40").Value = "'-2.45%"

$ws.Range("E41").Value = "'-6.35%"

$ws.Range("D42").Value = "'0.007080"
$ws.Range("E42").Value = "'8.04%"

$ws.Range("D43").Value = "'0.001953"
$ws.Range("E43").Value = "'-6.44%"

$ws.Range("D44").Value = "'0.008376"
$ws.Range("E44").Value = "'1.69%"

$ws.Range("D45").Value = "'0.3317"
$ws.Range("E45").Value = "'-3.85%"

$ws.Range("D46").Value = "'0.00006436"
$ws.Range("E46").Value = "'-4.69%"

$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'-0.07%"

$ws.Range("D49").Value = "'0.006331"
$ws.Range("E49").Value = "'86.43%"

$ws.Range("D50").Value = "'0.00002105"
$ws.Range("E50").Value = "'-0.07%"

$ws.Range("D51").Value = "'0.0002005"
$ws.Range("E51").Value = "'-0.07%"
